# Auto-generated script to apply scheduled price-runner updates
# to the Odin_Profits workbook (8 crafting-job sheets).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 660.2069
$ws.Cells.Item(12, 9).Value = 382.47827
$ws.Cells.Item(12, 10).Value = 1724.8334
$ws.Cells.Item(12, 11).Value = 382.47827
$ws.Cells.Item(12, 12).Value = 1724.8334
$ws.Cells.Item(12, 13).Value = -212.47827
$ws.Cells.Item(12, 14).Value = -2064.8334
$ws.Cells.Item(62, 8).Value = 7581333.5
$ws.Cells.Item(62, 9).Value = 8777598
$ws.Cells.Item(62, 10).Value = 4990
$ws.Cells.Item(62, 11).Value = 8777598
$ws.Cells.Item(62, 12).Value = 4990
$ws.Cells.Item(62, 13).Value = -8776974
$ws.Cells.Item(62, 14).Value = -6238
$ws.Cells.Item(65, 8).Value = 7581333.5
$ws.Cells.Item(65, 9).Value = 8777598
$ws.Cells.Item(65, 10).Value = 4990
$ws.Cells.Item(65, 11).Value = 43887990
$ws.Cells.Item(65, 12).Value = 24950
$ws.Cells.Item(65, 13).Value = -43884870
$ws.Cells.Item(65, 14).Value = -31190
$ws.Cells.Item(107, 8).Value = 2816.353
$ws.Cells.Item(107, 9).Value = 2993.2666
$ws.Cells.Item(107, 10).Value = 1489.5
$ws.Cells.Item(107, 11).Value = 2993.2666
$ws.Cells.Item(107, 12).Value = 1489.5
$ws.Cells.Item(107, 13).Value = -1073.2666
$ws.Cells.Item(107, 14).Value = -5329.5
$ws.Cells.Item(116, 8).Value = 15877587
$ws.Cells.Item(116, 9).Value = 27780152
$ws.Cells.Item(116, 10).Value = 7500
$ws.Cells.Item(116, 11).Value = 27780152
$ws.Cells.Item(116, 12).Value = 7500
$ws.Cells.Item(116, 13).Value = -27776710
$ws.Cells.Item(116, 14).Value = -14384
$ws.Cells.Item(132, 8).Value = 313265.12
$ws.Cells.Item(132, 9).Value = 335321.16
$ws.Cells.Item(132, 10).Value = 11832.333
$ws.Cells.Item(132, 11).Value = 1005963.48
$ws.Cells.Item(132, 12).Value = 35496.999
$ws.Cells.Item(132, 13).Value = -1003433.48
$ws.Cells.Item(132, 14).Value = -40556.999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 472.75
$ws.Cells.Item(4, 9).Value = 492.3
$ws.Cells.Item(4, 11).Value = 492.3
$ws.Cells.Item(4, 13).Value = -376.3
$ws.Cells.Item(88, 8).Value = 3885.3333
$ws.Cells.Item(88, 9).Value = 2237.4
$ws.Cells.Item(88, 11).Value = 2237.4
$ws.Cells.Item(88, 13).Value = -1831.4
$ws.Cells.Item(91, 8).Value = 3885.3333
$ws.Cells.Item(91, 9).Value = 2237.4
$ws.Cells.Item(91, 11).Value = 2237.4
$ws.Cells.Item(91, 13).Value = -833.4000000000001
$ws.Cells.Item(102, 8).Value = 2764.9656
$ws.Cells.Item(102, 9).Value = 2738.7144
$ws.Cells.Item(102, 11).Value = 2738.7144
$ws.Cells.Item(102, 13).Value = -1116.7144
$ws.Cells.Item(122, 8).Value = 2545.2942
$ws.Cells.Item(122, 9).Value = 2109.8845
$ws.Cells.Item(122, 11).Value = 6329.6535
$ws.Cells.Item(122, 13).Value = -3879.6535

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 19623236
$ws.Cells.Item(80, 9).Value = 1449.25
$ws.Cells.Item(80, 10).Value = 37064828
$ws.Cells.Item(80, 11).Value = 1449.25
$ws.Cells.Item(80, 12).Value = 37064828
$ws.Cells.Item(80, 13).Value = -451.25
$ws.Cells.Item(80, 14).Value = -37066824
$ws.Cells.Item(83, 8).Value = 19623236
$ws.Cells.Item(83, 9).Value = 1449.25
$ws.Cells.Item(83, 10).Value = 37064828
$ws.Cells.Item(83, 11).Value = 7246.25
$ws.Cells.Item(83, 12).Value = 185324140
$ws.Cells.Item(83, 13).Value = -2254.25
$ws.Cells.Item(83, 14).Value = -185334124
$ws.Cells.Item(86, 8).Value = 4457.0386
$ws.Cells.Item(86, 9).Value = 2050.4
$ws.Cells.Item(86, 10).Value = 7738.8184
$ws.Cells.Item(86, 11).Value = 2050.4
$ws.Cells.Item(86, 12).Value = 7738.8184
$ws.Cells.Item(86, 13).Value = -927.4000000000001
$ws.Cells.Item(86, 14).Value = -9984.8184
$ws.Cells.Item(89, 8).Value = 4457.0386
$ws.Cells.Item(89, 9).Value = 2050.4
$ws.Cells.Item(89, 10).Value = 7738.8184
$ws.Cells.Item(89, 11).Value = 10252
$ws.Cells.Item(89, 12).Value = 38694.092
$ws.Cells.Item(89, 13).Value = -4636
$ws.Cells.Item(89, 14).Value = -49926.092
$ws.Cells.Item(94, 8).Value = 29922.79
$ws.Cells.Item(94, 9).Value = 1823.2727
$ws.Cells.Item(94, 10).Value = 68559.625
$ws.Cells.Item(94, 11).Value = 1823.2727
$ws.Cells.Item(94, 12).Value = 68559.625
$ws.Cells.Item(94, 13).Value = -1372.2727
$ws.Cells.Item(94, 14).Value = -69461.625
$ws.Cells.Item(107, 8).Value = 2045250.9
$ws.Cells.Item(107, 9).Value = 2567802.8
$ws.Cells.Item(107, 10).Value = 7297.9
$ws.Cells.Item(107, 11).Value = 2567802.8
$ws.Cells.Item(107, 12).Value = 7297.9
$ws.Cells.Item(107, 13).Value = -2565882.8
$ws.Cells.Item(107, 14).Value = -11137.9
$ws.Cells.Item(134, 8).Value = 1259602.8
$ws.Cells.Item(134, 10).Value = 9286.777
$ws.Cells.Item(134, 12).Value = 27860.331
$ws.Cells.Item(134, 14).Value = -32930.331

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 9032.641
$ws.Cells.Item(31, 10).Value = 4393.6
$ws.Cells.Item(31, 12).Value = 4393.6
$ws.Cells.Item(31, 14).Value = -4983.6
$ws.Cells.Item(34, 8).Value = 9032.641
$ws.Cells.Item(34, 10).Value = 4393.6
$ws.Cells.Item(34, 12).Value = 4393.6
$ws.Cells.Item(34, 14).Value = -4797.6
$ws.Cells.Item(122, 8).Value = 5657.1904
$ws.Cells.Item(122, 9).Value = 2817.8823
$ws.Cells.Item(122, 11).Value = 8453.6469
$ws.Cells.Item(122, 13).Value = -6003.6469
$ws.Cells.Item(134, 8).Value = 190499310
$ws.Cells.Item(134, 9).Value = 380969060
$ws.Cells.Item(134, 10).Value = 29583.334
$ws.Cells.Item(134, 11).Value = 1142907180
$ws.Cells.Item(134, 12).Value = 88750.00199999999
$ws.Cells.Item(134, 13).Value = -1142904645
$ws.Cells.Item(134, 14).Value = -93820.00199999999
$ws.Cells.Item(141, 8).Value = 720000
$ws.Cells.Item(141, 10).Value = 990000
$ws.Cells.Item(141, 12).Value = 990000
$ws.Cells.Item(141, 14).Value = -1000360

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 7861088.5
$ws.Cells.Item(4, 9).Value = 7861088.5
$ws.Cells.Item(4, 11).Value = 23583265.5
$ws.Cells.Item(4, 13).Value = -23583153.5
$ws.Cells.Item(11, 8).Value = 189.6
$ws.Cells.Item(11, 9).Value = 174.5
$ws.Cells.Item(11, 10).Value = 250
$ws.Cells.Item(11, 11).Value = 523.5
$ws.Cells.Item(11, 12).Value = 750
$ws.Cells.Item(11, 13).Value = -383.5
$ws.Cells.Item(11, 14).Value = -1030
$ws.Cells.Item(86, 8).Value = 1760.8572
$ws.Cells.Item(86, 9).Value = 5095
$ws.Cells.Item(86, 10).Value = 427.2
$ws.Cells.Item(86, 11).Value = 15285
$ws.Cells.Item(86, 12).Value = 1281.6
$ws.Cells.Item(86, 13).Value = -14099
$ws.Cells.Item(86, 14).Value = -3653.6
$ws.Cells.Item(89, 8).Value = 1760.8572
$ws.Cells.Item(89, 9).Value = 5095
$ws.Cells.Item(89, 10).Value = 427.2
$ws.Cells.Item(89, 11).Value = 45855
$ws.Cells.Item(89, 12).Value = 3844.8
$ws.Cells.Item(89, 13).Value = -39927
$ws.Cells.Item(89, 14).Value = -15700.8
$ws.Cells.Item(140, 8).Value = 21741470
$ws.Cells.Item(140, 9).Value = 26317464
$ws.Cells.Item(140, 10).Value = 5497
$ws.Cells.Item(140, 11).Value = 78952392
$ws.Cells.Item(140, 12).Value = 16491
$ws.Cells.Item(140, 13).Value = -78947212
$ws.Cells.Item(140, 14).Value = -26851

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 741950.9
$ws.Cells.Item(102, 9).Value = 1187201
$ws.Cells.Item(102, 11).Value = 1187201
$ws.Cells.Item(102, 13).Value = -1185579
$ws.Cells.Item(122, 8).Value = 5049.3335
$ws.Cells.Item(122, 9).Value = 2941.7878
$ws.Cells.Item(122, 11).Value = 8825.3634
$ws.Cells.Item(122, 13).Value = -6375.3634
$ws.Cells.Item(126, 8).Value = 25008858
$ws.Cells.Item(126, 9).Value = 41670524
$ws.Cells.Item(126, 10).Value = 16363
$ws.Cells.Item(126, 11).Value = 125011572
$ws.Cells.Item(126, 12).Value = 49089
$ws.Cells.Item(126, 13).Value = -125009102
$ws.Cells.Item(126, 14).Value = -54029

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3820.182
$ws.Cells.Item(40, 9).Value = 3573.8572
$ws.Cells.Item(40, 10).Value = 4251.25
$ws.Cells.Item(40, 11).Value = 3573.8572
$ws.Cells.Item(40, 12).Value = 4251.25
$ws.Cells.Item(40, 13).Value = -3437.8572
$ws.Cells.Item(40, 14).Value = -4523.25
$ws.Cells.Item(93, 8).Value = 2610.45
$ws.Cells.Item(93, 9).Value = 2877.4614
$ws.Cells.Item(93, 10).Value = 2114.5715
$ws.Cells.Item(93, 11).Value = 2877.4614
$ws.Cells.Item(93, 12).Value = 2114.5715
$ws.Cells.Item(93, 13).Value = -1629.4614
$ws.Cells.Item(93, 14).Value = -4610.5715

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 14761.333
$ws.Cells.Item(62, 9).Value = 11149.8
$ws.Cells.Item(62, 10).Value = 19275.75
$ws.Cells.Item(62, 11).Value = 11149.8
$ws.Cells.Item(62, 12).Value = 19275.75
$ws.Cells.Item(62, 13).Value = -10525.8
$ws.Cells.Item(62, 14).Value = -20523.75
$ws.Cells.Item(65, 8).Value = 14761.333
$ws.Cells.Item(65, 9).Value = 11149.8
$ws.Cells.Item(65, 10).Value = 19275.75
$ws.Cells.Item(65, 11).Value = 55749
$ws.Cells.Item(65, 12).Value = 96378.75
$ws.Cells.Item(65, 13).Value = -52629
$ws.Cells.Item(65, 14).Value = -102618.75
$ws.Cells.Item(123, 8).Value = 57179
$ws.Cells.Item(123, 10).Value = 57179
$ws.Cells.Item(123, 12).Value = 57179
$ws.Cells.Item(123, 14).Value = -66979
$ws.Cells.Item(136, 8).Value = 17255880
$ws.Cells.Item(136, 9).Value = 35734140
$ws.Cells.Item(136, 11).Value = 107202420
$ws.Cells.Item(136, 13).Value = -107199870
